$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September log entry ("axis" @ 2024-09-18 08:12:44) was recorded,
# pushing the whole September_Details/September_Date history (and every
# row below it, including the trailing "Broadband" marker row) down by
# one row. Inserting a fresh row at row 40 reproduces that shift exactly,
# then we populate the new entry itself.
$ws.Rows.Item(40).Insert()

$ws.Range("R40").Value = "axis"
$ws.Range("S40").Value = "2024-09-18 08:12:44"
